# Updated priority list to have UC short names
#
# Rewrites each "N. <old description>" paragraph under the "Priority List
# of the functionalities" heading into "N. <UCShortName>", marking the new
# short-name token as a single spell-checked word (wrapped in
# <w:proofErr spellStart/.../spellEnd/>), matching what Word's proofing
# pass does for an unrecognised identifier. The "_GoBack" bookmark, which
# used to sit around "Link to Home view" (item 12), moves to surround the
# new "UploadSyllabus" token on item 1, since that is now the last place
# in the story touched by the edit.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"

function New-ParaXml {
    param(
        [string]$ParaId,
        [string]$Indent,
        [string]$Number,
        [string]$ShortName,
        [bool]$IncludeGoBack
    )

    $xml = "<w:p xmlns:w='$wNs' xmlns:w14='$w14Ns' w14:paraId='$ParaId' w14:textId='77777777' w:rsidR='00B84500' w:rsidRDefault='00B84500' w:rsidP='00B84500'>"
    $xml += "<w:pPr><w:ind w:left='$Indent'/><w:jc w:val='both'/></w:pPr>"
    $xml += "<w:r><w:t xml:space=`"preserve`">$Number. </w:t></w:r>"
    $xml += "<w:proofErr w:type='spellStart'/>"
    $xml += "<w:r><w:t>$ShortName</w:t></w:r>"
    if ($IncludeGoBack) {
        $xml += "<w:bookmarkStart w:id='14' w:name='_GoBack'/><w:bookmarkEnd w:id='14'/>"
    }
    $xml += "<w:proofErr w:type='spellEnd'/>"
    $xml += "</w:p>"
    return $xml
}

# (paragraph index in the current document, paraId, indent, number, old text
# (for a sanity check), short name, carries _GoBack)
$items = @(
    @{ Index = 51; ParaId = "3AF3137C"; Indent = "720"; Number = "5";  Old = "5. Ability to search syllabi";                            Short = "SearchExisting";      GoBack = $false },
    @{ Index = 52; ParaId = "1A3F5FA1"; Indent = "720"; Number = "8";  Old = "8. View syllabi details";                                  Short = "ViewSyllabusDetails"; GoBack = $false },
    @{ Index = 53; ParaId = "0C446654"; Indent = "720"; Number = "7";  Old = "7. Ability to download syllabi";                           Short = "DownloadSyllabus";    GoBack = $false },
    @{ Index = 54; ParaId = "371658E4"; Indent = "720"; Number = "2";  Old = "2. Ability to import syllabi for newly created courses";    Short = "ImportSyllabus";      GoBack = $false },
    @{ Index = 55; ParaId = "6F5CFB64"; Indent = "630"; Number = "10"; Old = "10. Administrator log on";                                 Short = "LogOn";               GoBack = $false },
    @{ Index = 56; ParaId = "5D251FFD"; Indent = "630"; Number = "11"; Old = "11. Administrator log off";                                Short = "LogOff";              GoBack = $false },
    @{ Index = 57; ParaId = "74492611"; Indent = "720"; Number = "3";  Old = "3. Edit existing syllabi";                                 Short = "EditSyllabus";        GoBack = $false },
    @{ Index = 58; ParaId = "35464F5D"; Indent = "720"; Number = "9";  Old = "9. Toggle active status of syllabi";                       Short = "ToggleActive";        GoBack = $false },
    @{ Index = 59; ParaId = "655F7FEB"; Indent = "720"; Number = "6";  Old = "6. Sort search results";                                   Short = "SortExisting";        GoBack = $false },
    @{ Index = 60; ParaId = "1BBB521C"; Indent = "720"; Number = "1";  Old = "1. Ability to upload syllabi for historical courses";       Short = "UploadSyllabus";      GoBack = $true  },
    @{ Index = 61; ParaId = "7C49B9AF"; Indent = "720"; Number = "4";  Old = "4. Delete existing syllabi";                               Short = "DeleteSyllabus";      GoBack = $false },
    @{ Index = 62; ParaId = "336573A6"; Indent = "630"; Number = "12"; Old = "12. Link to Home view";                                    Short = "ClickHomeLink";       GoBack = $false }
)

foreach ($item in $items) {
    $p = $d.Paragraphs($item.Index)
    # Range.Text includes the trailing paragraph mark (CR, char 13); strip it
    # before comparing against the plain-text sanity check.
    $actualText = $p.Range.Text.TrimEnd([char]13)
    if ($actualText -ne $item.Old) {
        throw "Paragraph $($item.Index): expected [$($item.Old)] but found [$actualText]"
    }
    $xml = New-ParaXml $item.ParaId $item.Indent $item.Number $item.Short $item.GoBack
    [void]$p.Range.InsertXML($xml)
}
